$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 795
$ws.Range("I2").Value = 1400
$ws.Range("J2").Value = 190
$ws.Range("K2").Value = 1400
$ws.Range("L2").Value = 190
$ws.Range("M2").Value = -1287
$ws.Range("N2").Value = -416
$ws.Range("H19").Value = 680.2222
$ws.Range("I19").Value = 691.5
$ws.Range("K19").Value = 691.5
$ws.Range("M19").Value = -516.5
$ws.Range("H33").Value = 510.85715
$ws.Range("I33").Value = 470.2
$ws.Range("K33").Value = 470.2
$ws.Range("M33").Value = -241.2
$ws.Range("H58").Value = 3569.7144
$ws.Range("I58").Value = 372
$ws.Range("J58").Value = 7833.3335
$ws.Range("K58").Value = 1116
$ws.Range("L58").Value = 23500.0005
$ws.Range("M58").Value = -966
$ws.Range("N58").Value = -23800.0005
$ws.Range("H88").Value = 155605.53
$ws.Range("I88").Value = 400878.6
$ws.Range("J88").Value = 2309.875
$ws.Range("K88").Value = 400878.6
$ws.Range("L88").Value = 2309.875
$ws.Range("M88").Value = -400472.6
$ws.Range("N88").Value = -3121.875
$ws.Range("H91").Value = 155605.53
$ws.Range("I91").Value = 400878.6
$ws.Range("J91").Value = 2309.875
$ws.Range("K91").Value = 400878.6
$ws.Range("L91").Value = 2309.875
$ws.Range("M91").Value = -399474.6
$ws.Range("N91").Value = -5117.875
$ws.Range("H92").Value = 2210.8667
$ws.Range("I92").Value = 820.06665
$ws.Range("J92").Value = 3601.6667
$ws.Range("K92").Value = 820.06665
$ws.Range("L92").Value = 3601.6667
$ws.Range("M92").Value = 427.93335
$ws.Range("N92").Value = -6097.6667
$ws.Range("H94").Value = 19131.625
$ws.Range("I94").Value = 3264.5
$ws.Range("J94").Value = 34998.75
$ws.Range("K94").Value = 3264.5
$ws.Range("L94").Value = 34998.75
$ws.Range("M94").Value = -2813.5
$ws.Range("N94").Value = -35900.75
$ws.Range("H97").Value = 1199.6
$ws.Range("J97").Value = 1199.6
$ws.Range("L97").Value = 3598.8
$ws.Range("N97").Value = -4590.799999999999
$ws.Range("H100").Value = 635
$ws.Range("I100").Value = 682.4545000000001
$ws.Range("K100").Value = 682.4545000000001
$ws.Range("M100").Value = -141.4545000000001
$ws.Range("H106").Value = 3970.7144
$ws.Range("I106").Value = 4049.1667
$ws.Range("K106").Value = 4049.1667
$ws.Range("M106").Value = -3418.1667
$ws.Range("H111").Value = 1004
$ws.Range("I111").Value = 876.6667
$ws.Range("J111").Value = 1195
$ws.Range("K111").Value = 2630.0001
$ws.Range("L111").Value = 3585
$ws.Range("M111").Value = 436.9998999999998
$ws.Range("N111").Value = -9719
$ws.Range("H112").Value = 3201.3572
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 3370.6924
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 10112.0772
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -12328.0772
$ws.Range("H113").Value = 6961.4443
$ws.Range("I113").Value = 3113.1667
$ws.Range("J113").Value = 14658
$ws.Range("K113").Value = 3113.1667
$ws.Range("L113").Value = 14658
$ws.Range("M113").Value = 140.8332999999998
$ws.Range("N113").Value = -21166
$ws.Range("H127").Value = 2293.4
$ws.Range("I127").Value = 2293.4
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 6880.200000000001
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -1920.200000000001
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 222860.95
$ws.Range("I132").Value = 605.0769
$ws.Range("K132").Value = 1815.2307
$ws.Range("M132").Value = 714.7692999999999
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 7790.857
$ws.Range("I135").Value = 697.9231
$ws.Range("K135").Value = 6281.3079
$ws.Range("M135").Value = -3746.3079
$ws.Range("H137").Value = 6842.7144
$ws.Range("I137").Value = 4533.5
$ws.Range("K137").Value = 13600.5
$ws.Range("M137").Value = -11050.5
$ws.Range("H138").Value = 3361.6724
$ws.Range("I138").Value = 1580.0938
$ws.Range("J138").Value = 5554.385
$ws.Range("K138").Value = 4740.2814
$ws.Range("L138").Value = 16663.155
$ws.Range("M138").Value = 399.7186000000002
$ws.Range("N138").Value = -26943.155
$ws.Range("H141").Value = 4952.8213
$ws.Range("I141").Value = 4702.4707
$ws.Range("K141").Value = 14107.4121
$ws.Range("M141").Value = -8927.4121

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H5").Value = 1316.3334
$ws.Range("I5").Value = 999.5
$ws.Range("J5").Value = 1950
$ws.Range("K5").Value = 999.5
$ws.Range("L5").Value = 1950
$ws.Range("M5").Value = -887.5
$ws.Range("N5").Value = -2174
$ws.Range("H21").Value = 1806.2858
$ws.Range("I21").Value = 928.8
$ws.Range("J21").Value = 4000
$ws.Range("K21").Value = 928.8
$ws.Range("L21").Value = 4000
$ws.Range("M21").Value = -554.8
$ws.Range("N21").Value = -4748
$ws.Range("H32").Value = 4407.8423
$ws.Range("I32").Value = 2152.889
$ws.Range("K32").Value = 2152.889
$ws.Range("M32").Value = -1865.889
$ws.Range("H45").Value = 3523.8
$ws.Range("I45").Value = 5337.3335
$ws.Range("J45").Value = 803.5
$ws.Range("K45").Value = 5337.3335
$ws.Range("L45").Value = 803.5
$ws.Range("M45").Value = -4960.3335
$ws.Range("N45").Value = -1557.5
$ws.Range("H61").Value = 2833.5
$ws.Range("I61").Value = 2380.7334
$ws.Range("K61").Value = 2380.7334
$ws.Range("M61").Value = -2168.7334
$ws.Range("H74").Value = 867.9375
$ws.Range("I74").Value = 828.2727
$ws.Range("K74").Value = 828.2727
$ws.Range("M74").Value = 45.72730000000001
$ws.Range("H77").Value = 867.9375
$ws.Range("I77").Value = 828.2727
$ws.Range("K77").Value = 4141.363499999999
$ws.Range("M77").Value = 226.6365000000005
$ws.Range("H88").Value = 2336.0833
$ws.Range("I88").Value = 1817.6
$ws.Range("J88").Value = 2706.4285
$ws.Range("K88").Value = 1817.6
$ws.Range("L88").Value = 2706.4285
$ws.Range("M88").Value = -1411.6
$ws.Range("N88").Value = -3518.4285
$ws.Range("H91").Value = 2336.0833
$ws.Range("I91").Value = 1817.6
$ws.Range("J91").Value = 2706.4285
$ws.Range("K91").Value = 1817.6
$ws.Range("L91").Value = 2706.4285
$ws.Range("M91").Value = -413.5999999999999
$ws.Range("N91").Value = -5514.4285
$ws.Range("H97").Value = 2251.4285
$ws.Range("I97").Value = 1580.875
$ws.Range("J97").Value = 4397.2
$ws.Range("K97").Value = 1580.875
$ws.Range("L97").Value = 4397.2
$ws.Range("M97").Value = -1084.875
$ws.Range("N97").Value = -5389.2
$ws.Range("H109").Value = 65141.715
$ws.Range("J109").Value = 65141.715
$ws.Range("L109").Value = 65141.715
$ws.Range("N109").Value = -67915.715
$ws.Range("H110").Value = 1433.7273
$ws.Range("I110").Value = 1284.9524
$ws.Range("K110").Value = 1284.9524
$ws.Range("M110").Value = 760.0476000000001
$ws.Range("H132").Value = 1754.5
$ws.Range("I132").Value = 1717.0968
$ws.Range("J132").Value = 2914
$ws.Range("K132").Value = 5151.2904
$ws.Range("L132").Value = 8742
$ws.Range("M132").Value = -2621.2904
$ws.Range("N132").Value = -13802
$ws.Range("H136").Value = 2833.5
$ws.Range("I136").Value = 2380.7334
$ws.Range("K136").Value = 7142.2002
$ws.Range("M136").Value = -4592.2002

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H4").Value = 1316.3334
$ws.Range("I4").Value = 999.5
$ws.Range("J4").Value = 1950
$ws.Range("K4").Value = 999.5
$ws.Range("L4").Value = 1950
$ws.Range("M4").Value = -884.5
$ws.Range("N4").Value = -2180
$ws.Range("H94").Value = 1930.7646
$ws.Range("I94").Value = 2374
$ws.Range("J94").Value = 1432.125
$ws.Range("K94").Value = 2374
$ws.Range("L94").Value = 1432.125
$ws.Range("M94").Value = -1923
$ws.Range("N94").Value = -2334.125
$ws.Range("H105").Value = 3163.2307
$ws.Range("I105").Value = 2527.2778
$ws.Range("K105").Value = 2527.2778
$ws.Range("M105").Value = -780.2777999999998
$ws.Range("H134").Value = 2347.3547
$ws.Range("I134").Value = 2209.2415
$ws.Range("K134").Value = 6627.7245
$ws.Range("M134").Value = -4092.7245

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H4").Value = 31753.77
$ws.Range("I4").Value = 14600
$ws.Range("J4").Value = 42474.875
$ws.Range("K4").Value = 14600
$ws.Range("L4").Value = 42474.875
$ws.Range("M4").Value = -14488
$ws.Range("N4").Value = -42698.875
$ws.Range("H7").Value = 1125
$ws.Range("I7").Value = 1233.3334
$ws.Range("J7").Value = 800
$ws.Range("K7").Value = 1233.3334
$ws.Range("L7").Value = 800
$ws.Range("M7").Value = -1120.3334
$ws.Range("N7").Value = -1026
$ws.Range("H16").Value = 1267.3636
$ws.Range("I16").Value = 1215.6666
$ws.Range("K16").Value = 1215.6666
$ws.Range("M16").Value = -928.6666
$ws.Range("H31").Value = 6172.5864
$ws.Range("I31").Value = 4013.2
$ws.Range("K31").Value = 4013.2
$ws.Range("M31").Value = -3718.2
$ws.Range("H34").Value = 6172.5864
$ws.Range("I34").Value = 4013.2
$ws.Range("K34").Value = 4013.2
$ws.Range("M34").Value = -3811.2
$ws.Range("H52").Value = 39999
$ws.Range("J52").Value = 39999
$ws.Range("L52").Value = 39999
$ws.Range("N52").Value = -40587
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H94").Value = 778
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 778
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 778
$ws.Range("N94").Value = -1680
$ws.Range("M94").ClearContents()
$ws.Range("H105").Value = 1456.1666
$ws.Range("I105").Value = 1297.9286
$ws.Range("J105").Value = 2010
$ws.Range("K105").Value = 1297.9286
$ws.Range("L105").Value = 2010
$ws.Range("M105").Value = 449.0714
$ws.Range("N105").Value = -5504
$ws.Range("H113").Value = 1267.3636
$ws.Range("I113").Value = 1215.6666
$ws.Range("K113").Value = 1215.6666
$ws.Range("M113").Value = 954.3334
$ws.Range("H122").Value = 4193.1577
$ws.Range("I122").Value = 3992.353
$ws.Range("J122").Value = 5900
$ws.Range("K122").Value = 11977.059
$ws.Range("L122").Value = 17700
$ws.Range("M122").Value = -9527.059000000001
$ws.Range("N122").Value = -22600
$ws.Range("H132").Value = 1182.6765
$ws.Range("I132").Value = 976.9286
$ws.Range("K132").Value = 2930.7858
$ws.Range("M132").Value = -400.7857999999997
$ws.Range("H134").Value = 1603.8214
$ws.Range("I134").Value = 1533.5927
$ws.Range("K134").Value = 4600.7781
$ws.Range("M134").Value = -2065.7781
$ws.Range("H139").Value = 144975
$ws.Range("J139").Value = 144975
$ws.Range("L139").Value = 144975
$ws.Range("N139").Value = -155255
$ws.Range("H141").Value = 265148
$ws.Range("I141").Value = 30296
$ws.Range("J141").Value = 500000
$ws.Range("K141").Value = 30296
$ws.Range("L141").Value = 500000
$ws.Range("M141").Value = -25116
$ws.Range("N141").Value = -510360

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 1177.5
$ws.Range("J5").Value = 2005
$ws.Range("L5").Value = 6015
$ws.Range("N5").Value = -6239
$ws.Range("H8").Value = 495.33334
$ws.Range("I8").Value = 495.33334
$ws.Range("K8").Value = 1486.00002
$ws.Range("M8").Value = -1347.00002
$ws.Range("H12").Value = 768
$ws.Range("J12").Value = 866.3333
$ws.Range("L12").Value = 2598.9999
$ws.Range("N12").Value = -2944.9999
$ws.Range("H45").Value = 4567.6665
$ws.Range("J45").Value = 4537
$ws.Range("L45").Value = 13611
$ws.Range("N45").Value = -14675
$ws.Range("H81").Value = 4499.3335
$ws.Range("I81").Value = 4300
$ws.Range("J81").Value = 4599
$ws.Range("K81").Value = 12900
$ws.Range("L81").Value = 13797
$ws.Range("M81").Value = -11777
$ws.Range("N81").Value = -16043
$ws.Range("H84").Value = 4499.3335
$ws.Range("I84").Value = 4300
$ws.Range("J84").Value = 4599
$ws.Range("K84").Value = 38700
$ws.Range("L84").Value = 41391
$ws.Range("M84").Value = -33084
$ws.Range("N84").Value = -52623
$ws.Range("H107").Value = 1057.8
$ws.Range("I107").Value = 645
$ws.Range("J107").Value = 1333
$ws.Range("K107").Value = 1935
$ws.Range("L107").Value = 3999
$ws.Range("M107").Value = -15
$ws.Range("N107").Value = -7839
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H116").Value = 100000
$ws.Range("I116").Value = 100000
$ws.Range("K116").Value = 300000
$ws.Range("M116").Value = -296558
$ws.Range("H117").Value = 399
$ws.Range("I117").Value = 399
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 1197
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 2245
$ws.Range("N117").ClearContents()
$ws.Range("H120").Value = 3706.5
$ws.Range("I120").Value = 3706.5
$ws.Range("K120").Value = 11119.5
$ws.Range("M120").Value = -6281.5
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H131").Value = 4223872
$ws.Range("I131").Value = 223160.2
$ws.Range("K131").Value = 669480.6000000001
$ws.Range("M131").Value = -664440.6000000001
$ws.Range("H133").Value = 5998.5
$ws.Range("I133").Value = 5998.5
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 17995.5
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -12935.5
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 1177.5
$ws.Range("J135").Value = 2005
$ws.Range("L135").Value = 18045
$ws.Range("N135").Value = -23115
$ws.Range("H138").Value = 1646
$ws.Range("I138").Value = 1646
$ws.Range("K138").Value = 4938
$ws.Range("M138").Value = 202
$ws.Range("H139").Value = 4720.9287
$ws.Range("I139").Value = 1584.7142
$ws.Range("K139").Value = 4754.142599999999
$ws.Range("M139").Value = 385.8574000000008
$ws.Range("H141").Value = 1514
$ws.Range("I141").Value = 1514
$ws.Range("K141").Value = 4542
$ws.Range("M141").Value = 638

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 122.38461
$ws.Range("J2").Value = 50.8
$ws.Range("L2").Value = 50.8
$ws.Range("N2").Value = -276.8
$ws.Range("H33").Value = 20880.8
$ws.Range("I33").Value = 20017
$ws.Range("J33").Value = 21096.75
$ws.Range("K33").Value = 20017
$ws.Range("L33").Value = 21096.75
$ws.Range("M33").Value = -19765
$ws.Range("N33").Value = -21600.75
$ws.Range("H63").Value = 50114
$ws.Range("J63").Value = 50114
$ws.Range("L63").Value = 50114
$ws.Range("N63").Value = -51486
$ws.Range("H66").Value = 50114
$ws.Range("J66").Value = 50114
$ws.Range("L66").Value = 150342
$ws.Range("N66").Value = -157206
$ws.Range("H113").Value = 1433
$ws.Range("I113").Value = 1399.5
$ws.Range("K113").Value = 1399.5
$ws.Range("M113").Value = 770.5
$ws.Range("H122").Value = 1033.762
$ws.Range("I122").Value = 1018.55554
$ws.Range("K122").Value = 3055.66662
$ws.Range("M122").Value = -605.66662
$ws.Range("H126").Value = 3613.4285
$ws.Range("I126").Value = 2857.5
$ws.Range("J126").Value = 4621.3335
$ws.Range("K126").Value = 8572.5
$ws.Range("L126").Value = 13864.0005
$ws.Range("M126").Value = -6102.5
$ws.Range("N126").Value = -18804.0005
$ws.Range("H132").Value = 3343.625
$ws.Range("I132").Value = 3343.625
$ws.Range("K132").Value = 10030.875
$ws.Range("M132").Value = -7500.875

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 2443.5
$ws.Range("I16").Value = 610
$ws.Range("K16").Value = 610
$ws.Range("M16").Value = -440
$ws.Range("H22").Value = 1866
$ws.Range("I22").Value = 1279.8
$ws.Range("J22").Value = 3331.5
$ws.Range("K22").Value = 1279.8
$ws.Range("L22").Value = 3331.5
$ws.Range("M22").Value = -984.8
$ws.Range("N22").Value = -3921.5
$ws.Range("H27").Value = 1866
$ws.Range("I27").Value = 1279.8
$ws.Range("J27").Value = 3331.5
$ws.Range("K27").Value = 1279.8
$ws.Range("L27").Value = 3331.5
$ws.Range("M27").Value = -1172.8
$ws.Range("N27").Value = -3545.5
$ws.Range("H46").Value = 1073.8
$ws.Range("I46").Value = 947.3125
$ws.Range("J46").Value = 1298.6666
$ws.Range("K46").Value = 947.3125
$ws.Range("L46").Value = 1298.6666
$ws.Range("M46").Value = -759.3125
$ws.Range("N46").Value = -1674.6666
$ws.Range("H64").Value = 2075
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 2150
$ws.Range("K64").Value = 2000
$ws.Range("L64").Value = 2150
$ws.Range("M64").Value = -1775
$ws.Range("N64").Value = -2600
$ws.Range("H67").Value = 2075
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 2150
$ws.Range("K67").Value = 2000
$ws.Range("L67").Value = 2150
$ws.Range("M67").Value = -1220
$ws.Range("N67").Value = -3710
$ws.Range("H75").Value = 49999
$ws.Range("I75").Value = 49999
$ws.Range("K75").Value = 49999
$ws.Range("M75").Value = -49063
$ws.Range("H78").Value = 49999
$ws.Range("I78").Value = 49999
$ws.Range("K78").Value = 149997
$ws.Range("M78").Value = -145317
$ws.Range("H93").Value = 3221.8333
$ws.Range("I93").Value = 3566.2
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 3566.2
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = -2318.2
$ws.Range("N93").Value = -3996
$ws.Range("H116").Value = 50000
$ws.Range("J116").Value = 50000
$ws.Range("L116").Value = 50000
$ws.Range("N116").Value = -59178
$ws.Range("H122").Value = 2675.158
$ws.Range("I122").Value = 2414.6667
$ws.Range("K122").Value = 7244.000100000001
$ws.Range("M122").Value = -4794.000100000001
$ws.Range("H132").Value = 2035.2646
$ws.Range("I132").Value = 1712.9231
$ws.Range("K132").Value = 5138.7693
$ws.Range("M132").Value = -2608.7693
$ws.Range("H136").Value = 910.1667
$ws.Range("I136").Value = 910.1667
$ws.Range("K136").Value = 2730.5001
$ws.Range("M136").Value = -180.5001000000002

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 6524.091
$ws.Range("J122").Value = 9999.5
$ws.Range("L122").Value = 29998.5
$ws.Range("N122").Value = -34898.5
$ws.Range("H132").Value = 2205.6296
$ws.Range("I132").Value = 2232.8684
$ws.Range("J132").Value = 2140.9375
$ws.Range("K132").Value = 6698.6052
$ws.Range("L132").Value = 6422.8125
$ws.Range("M132").Value = -4168.6052
$ws.Range("N132").Value = -11482.8125
$ws.Range("H136").Value = 664.3077
$ws.Range("I136").Value = 605.9722
$ws.Range("K136").Value = 1817.9166
$ws.Range("M136").Value = 732.0834

